$p = $ppt.ActivePresentation

# --- Slide 1 (title slide): update the tutorial/venue line in the subtitle ---
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item(2)
$subtitleText = $subtitle.TextFrame.TextRange
$venuePara = $subtitleText.Paragraphs(2)
$venueRun = $venuePara.Runs(1)
$venueRun.Text = "Software Productivity and Sustainability track, ATPESC 2021"

# --- Slide 2 (license/citation slide): update citation text and DOI ---
$slide2 = $p.Slides.Item(2)
$content = $slide2.Shapes.Item(2)
$contentText = $content.TextFrame.TextRange
$citationPara = $contentText.Paragraphs(3)

$citationRun1 = $citationPara.Runs(1)
$citationRun1.Text = "The requested citation the overall tutorial is: David E. Bernholdt, Anshu Dubey, Rinku K. Gupta, and David M. Rogers, Software Productivity and Sustainability track, in Argonne Training Program on Extreme-Scale Computing (ATPESC), online, 2021. DOI: "

$citationRun2 = $citationPara.Runs(2)
$citationRun2.Text = "10.6084/m9.figshare.15130590"
